# 9th Stab- Cosmetic Changes
# A new "Jun_17" / "Jun_15" pair of rating-date columns is inserted in front of
# the existing "Jun_13" / "Jun_10" columns, shifting the report window to the
# right. All rows default to "UN" (unchanged) for the two newly-added columns,
# except for BidaskClub (row 22) which received a new rating change on
# 6/13/2018 that is highlighted with a fill color.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh columns before column B; this shifts the old "Jun_13"
# column (B) to D and the old "Jun_10" column (C) to E, preserving all of
# their existing values/styles.
$ws.Columns("B:B").Insert()
$ws.Columns("B:B").Insert()

# Keep the same custom column width (8.0) used by the original "Jun_10"
# column across the whole C:E block, now that it spans three columns.
$ws.Columns("C:C").ColumnWidth = 7.1666666666667
$ws.Columns("D:D").ColumnWidth = 7.1666666666667
$ws.Columns("E:E").ColumnWidth = 7.1666666666667

# --- Header row -----------------------------------------------------------
$ws.Cells.Item(1, 2).Value = "Jun_17"
$ws.Cells.Item(1, 3).Value = "Jun_15"

# --- Data rows --------------------------------------------------------------
# Default every new cell in columns B (Jun_17) and C (Jun_15) to "UN" (no
# analyst rating change reported for that date).
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# BidaskClub (row 22) got a new rating change reported on 6/13/2018 for the
# Jun_15 column; highlight it with a fill color, same as analysts usually
# get when a new rating change shows up.
$ws.Cells.Item(22, 3).Value = "6/13/2018,Upgrades,Strong Sell -> Sell,"
$ws.Cells.Item(22, 3).Interior.Pattern = 1
$ws.Cells.Item(22, 3).Interior.ColorIndex = 42

# Roth Capital (row 23) had initiated a rating on 6/13/2018 that is still in
# effect for the Jun_15 column as well, so it carries the same note forward
# (no highlight, since it isn't a brand new change for that column).
$ws.Cells.Item(23, 3).Value = "6/13/2018,Initiates,Buy,`$31.00"
